{"js": "// Remove the trailing Jekyll-site boilerplate paragraphs that used to follow\n// the bibliography entry ending in \"Editora Protec, 1991. PROVENZA, F.\n// Projetista de M\u00e1quinas . Editora Protec, 1991.\":\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"Editora Protec, 1991. PROVENZA, F. Projetista de M\u00e1quinas . Editora Protec, 1991.\";\n// The exact sequence of paragraph texts expected to immediately follow the\n// marker paragraph (in order). Only delete when the whole pattern matches,\n// so the script is a no-op if the boilerplate is already gone.\nconst expectedSequence = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  const start = markerIndex + 1;\n  let matches = start + expectedSequence.length <= items.length;\n  if (matches) {\n    for (let k = 0; k < expectedSequence.length; k++) {\n      if (items[start + k].text.trim() !== expectedSequence[k]) {\n        matches = false;\n        break;\n      }\n    }\n  }\n  if (matches) {\n    for (let k = 0; k < expectedSequence.length; k++) {\n      items[start + k].delete();\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing Jekyll-site boilerplate paragraphs that used to follow\n# the bibliography entry ending in \"Editora Protec, 1991. PROVENZA, F.\n# Projetista de M\u00e1quinas . Editora Protec, 1991.\":\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n\n$d = $word.ActiveDocument\n\n$marker = \"Editora Protec, 1991. PROVENZA, F. Projetista de M\u00e1quinas . Editora Protec, 1991.\"\n# The exact sequence of paragraph texts expected to immediately follow the\n# marker paragraph (in order). Only delete when the whole pattern matches,\n# so the script is a no-op if the boilerplate is already gone.\n$expectedSequence = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$count = $d.Paragraphs.Count\n$markerIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $marker) {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -ge 1) {\n    $start = $markerIndex + 1\n    $matches = ($start + $expectedSequence.Length - 1) -le $count\n    if ($matches) {\n        for ($k = 0; $k -lt $expectedSequence.Length; $k++) {\n            $t = $d.Paragraphs.Item($start + $k).Range.Text.TrimEnd(\"`r\", \"`a\")\n            if ($t -ne $expectedSequence[$k]) {\n                $matches = $false\n                break\n            }\n        }\n    }\n    if ($matches) {\n        for ($k = 0; $k -lt $expectedSequence.Length; $k++) {\n            # Paragraph at $start always addresses the next surviving\n            # paragraph of the block since earlier deletes shift it down.\n            $d.Paragraphs.Item($start).Range.Delete()\n        }\n    }\n}\n"}
